$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.97
$ws.Range("L2").Value = 1.25
$ws.Range("Y2").Value = 25
$ws.Range("AF2").Value = 16
$ws.Range("AJ2").Value = 25
$ws.Range("AN2").Value = 8.6
$ws.Range("N3").Value = 1.1
$ws.Range("V3").Value = 1.14
$ws.Range("F4").Value = 4.6
$ws.Range("G4").Value = 5.4
$ws.Range("H4").Value = 1.65
$ws.Range("I4").Value = 1.75
$ws.Range("K4").Value = 5
$ws.Range("N4").Value = 5.7
$ws.Range("F5").Value = 2.2
$ws.Range("I5").Value = 4.9
$ws.Range("J5").Value = 2.64
$ws.Range("K5").Value = 3.2
$ws.Range("N5").Value = 2.28
$ws.Range("O5").Value = 1.62
$ws.Range("S5").Value = 5.6
$ws.Range("X5").Value = 8.800000000000001
$ws.Range("AB5").Value = 14.5
$ws.Range("G6").Value = 1.48
$ws.Range("H6").Value = 8.800000000000001
$ws.Range("L6").Value = 1.37
$ws.Range("P6").Value = 1.69
$ws.Range("Q6").Value = 2.16
$ws.Range("T6").Value = 2.46
$ws.Range("W6").Value = 3.05
$ws.Range("X6").Value = 990
$ws.Range("AN6").Value = 12.5
$ws.Range("G7").Value = 1.5
$ws.Range("J7").Value = 1.09
$ws.Range("N7").Value = 1.1
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.07
$ws.Range("W7").Value = 1.01
$ws.Range("G8").Value = 1.4
$ws.Range("T8").Value = 1.81
$ws.Range("W8").Value = 3.45
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 1.78
$ws.Range("H10").Value = 3
$ws.Range("F11").Value = 2.04
$ws.Range("G11").Value = 2.12
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 3.85
$ws.Range("J11").Value = 4
$ws.Range("O11").Value = 1.18
$ws.Range("Q11").Value = 1.53
$ws.Range("S11").Value = 2.32
$ws.Range("U11").Value = 2.52
$ws.Range("V11").Value = 1.36
$ws.Range("W11").Value = 1.89
$ws.Range("F12").Value = 2.38
$ws.Range("G12").Value = 2.48
$ws.Range("H12").Value = 3.15
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.6
$ws.Range("P12").Value = 1.73
$ws.Range("Q12").Value = 1.96
$ws.Range("W12").Value = 1.67
$ws.Range("X12").Value = 14.5
$ws.Range("AC12").Value = 9.199999999999999
$ws.Range("AG12").Value = 14.5
$ws.Range("H13").Value = 1.44
$ws.Range("Q13").Value = 1.57
$ws.Range("X13").Value = 30
$ws.Range("AO13").Value = 6.6
$ws.Range("F14").Value = 2.2
$ws.Range("G14").Value = 2.36
$ws.Range("H14").Value = 3.35
$ws.Range("P14").Value = 1.99
$ws.Range("Q14").Value = 1.84
$ws.Range("W14").Value = 1.74
$ws.Range("AB14").Value = 13
$ws.Range("AC14").Value = 10
$ws.Range("I15").Value = 14.5
$ws.Range("J15").Value = 7.2
$ws.Range("N15").Value = 7
$ws.Range("AF15").Value = 980
$ws.Range("F16").Value = 2.16
$ws.Range("I16").Value = 3.85
$ws.Range("Q16").Value = 1.78
$ws.Range("R16").Value = 1.43
$ws.Range("AJ16").Value = 980
$ws.Range("AK16").Value = 27
$ws.Range("AL16").Value = 40
$ws.Range("AN16").Value = 17
$ws.Range("G17").Value = 6.4
$ws.Range("H17").Value = 1.63
$ws.Range("F18").Value = 2.6
$ws.Range("I18").Value = 3.1
$ws.Range("K18").Value = 3.6
$ws.Range("P18").Value = 1.81
$ws.Range("V18").Value = 1.48
$ws.Range("W18").Value = 1.55
$ws.Range("F19").Value = 4.3
$ws.Range("H19").Value = 1.75
$ws.Range("Q19").Value = 1.57
$ws.Range("R19").Value = 1.46
$ws.Range("S19").Value = 2.42
$ws.Range("U19").Value = 1.89
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 1.48
$ws.Range("I20").Value = 1.58
$ws.Range("K20").Value = 5.8
$ws.Range("M20").Value = 1.02
$ws.Range("N20").Value = 7
$ws.Range("P20").Value = 3
$ws.Range("R20").Value = 1.82
$ws.Range("S20").Value = 2.02
$ws.Range("T20").Value = 1.61
$ws.Range("Z20").Value = 980
$ws.Range("AL20").Value = 65
$ws.Range("H21").Value = 1.72
$ws.Range("P21").Value = 1.94
$ws.Range("Q21").Value = 1.88
$ws.Range("W21").Value = 1.21
$ws.Range("I22").Value = 3.75
$ws.Range("F23").Value = 2.22
$ws.Range("G23").Value = 2.42
$ws.Range("H23").Value = 2.96
$ws.Range("I23").Value = 3.35
$ws.Range("J23").Value = 3.85
$ws.Range("K23").Value = 4.4
$ws.Range("M23").Value = 1.02
$ws.Range("Q23").Value = 1.58
$ws.Range("S23").Value = 2.44
$ws.Range("T23").Value = 1.6
$ws.Range("V23").Value = 1.42
$ws.Range("W23").Value = 1.71
$ws.Range("AA23").Value = 60
$ws.Range("AC23").Value = 12
$ws.Range("AG23").Value = 14.5
$ws.Range("AK23").Value = 26
$ws.Range("F24").Value = 1.2
$ws.Range("S24").Value = 1.69
$ws.Range("T24").Value = 1.64
$ws.Range("W24").Value = 5.1
$ws.Range("X24").Value = 1000
$ws.Range("Y24").Value = 100
$ws.Range("Z24").Value = 210
$ws.Range("AA24").Value = 610
$ws.Range("AD24").Value = 50
$ws.Range("AE24").Value = 190
$ws.Range("AG24").Value = 13.5
$ws.Range("AH24").Value = 30
$ws.Range("AN24").Value = 2.86
$ws.Range("H25").Value = 2.76
$ws.Range("P25").Value = 2.42
$ws.Range("Q25").Value = 1.69
$ws.Range("T25").Value = 1.56
$ws.Range("U25").Value = 2.7
$ws.Range("X25").Value = 20
$ws.Range("AO25").Value = 17.5
$ws.Range("AD16").Value = 17.5
$ws.Range("AF16").Value = 18.5
